$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Cells.Item(70,8).Value = 12502939
$ws.Cells.Item(70,9).Value = 3164.6667
$ws.Cells.Item(70,10).Value = 14708782
$ws.Cells.Item(70,11).Value = 9494.000100000001
$ws.Cells.Item(70,12).Value = 44126346
$ws.Cells.Item(70,13).Value = -9224.000100000001
$ws.Cells.Item(70,14).Value = -44126886
# Row 73
$ws.Cells.Item(73,8).Value = 12502939
$ws.Cells.Item(73,9).Value = 3164.6667
$ws.Cells.Item(73,10).Value = 14708782
$ws.Cells.Item(73,11).Value = 9494.000100000001
$ws.Cells.Item(73,12).Value = 44126346
$ws.Cells.Item(73,13).Value = -8558.000100000001
$ws.Cells.Item(73,14).Value = -44128218
# Row 94
$ws.Cells.Item(94,8).Value = 426
$ws.Cells.Item(94,9).Value = 426
$ws.Cells.Item(94,11).Value = 426
$ws.Cells.Item(94,13).Value = 25
# Row 103
$ws.Cells.Item(103,8).Value = 1585.7142
$ws.Cells.Item(103,10).Value = 1350
$ws.Cells.Item(103,12).Value = 4050
$ws.Cells.Item(103,14).Value = -5222
# Row 125
$ws.Cells.Item(125,8).Value = 1217.0667
$ws.Cells.Item(125,9).Value = 1107.4
$ws.Cells.Item(125,10).Value = 1271.9
$ws.Cells.Item(125,11).Value = 9966.6
$ws.Cells.Item(125,12).Value = 11447.1
$ws.Cells.Item(125,13).Value = -7506.6
$ws.Cells.Item(125,14).Value = -16367.1
# Row 138
$ws.Cells.Item(138,8).Value = 3042.9272
$ws.Cells.Item(138,9).Value = 2569.1765
$ws.Cells.Item(138,11).Value = 7707.529500000001
$ws.Cells.Item(138,13).Value = -2567.529500000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61,8).Value = 2395.8125
$ws.Cells.Item(61,9).Value = 2104.1785
$ws.Cells.Item(61,11).Value = 2104.1785
$ws.Cells.Item(61,13).Value = -1892.1785
# Row 97
$ws.Cells.Item(97,8).Value = 1184.7037
$ws.Cells.Item(97,9).Value = 870.34784
$ws.Cells.Item(97,11).Value = 870.34784
$ws.Cells.Item(97,13).Value = -374.34784
# Row 136
$ws.Cells.Item(136,8).Value = 2395.8125
$ws.Cells.Item(136,9).Value = 2104.1785
$ws.Cells.Item(136,11).Value = 6312.5355
$ws.Cells.Item(136,13).Value = -3762.5355

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 60
$ws.Cells.Item(60,8).Value = 106945.664
$ws.Cells.Item(60,10).Value = 106945.664
$ws.Cells.Item(60,12).Value = 106945.664
$ws.Cells.Item(60,14).Value = -108143.664
# Row 82
$ws.Cells.Item(82,8).Value = 12097.4
$ws.Cells.Item(82,9).Value = 5124.25
$ws.Cells.Item(82,10).Value = 39990
$ws.Cells.Item(82,11).Value = 5124.25
$ws.Cells.Item(82,12).Value = 39990
$ws.Cells.Item(82,13).Value = -4741.25
$ws.Cells.Item(82,14).Value = -40756
# Row 85
$ws.Cells.Item(85,8).Value = 12097.4
$ws.Cells.Item(85,9).Value = 5124.25
$ws.Cells.Item(85,10).Value = 39990
$ws.Cells.Item(85,11).Value = 5124.25
$ws.Cells.Item(85,12).Value = 39990
$ws.Cells.Item(85,13).Value = -3798.25
$ws.Cells.Item(85,14).Value = -42642
# Row 105
$ws.Cells.Item(105,8).Value = 3971.6667
$ws.Cells.Item(105,9).Value = 1995
$ws.Cells.Item(105,10).Value = 4218.75
$ws.Cells.Item(105,11).Value = 1995
$ws.Cells.Item(105,12).Value = 4218.75
$ws.Cells.Item(105,13).Value = -248
$ws.Cells.Item(105,14).Value = -7712.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7,8).Value = 180.1875
$ws.Cells.Item(7,10).Value = 1000
$ws.Cells.Item(7,12).Value = 1000
$ws.Cells.Item(7,14).Value = -1226
# Row 107
$ws.Cells.Item(107,8).Value = 92742
$ws.Cells.Item(107,9).Value = 143523.14
$ws.Cells.Item(107,10).Value = 3875
$ws.Cells.Item(107,11).Value = 143523.14
$ws.Cells.Item(107,12).Value = 3875
$ws.Cells.Item(107,13).Value = -141603.14
$ws.Cells.Item(107,14).Value = -7715
# Row 122
$ws.Cells.Item(122,8).Value = 3038.1667
$ws.Cells.Item(122,9).Value = 3038.1667
$ws.Cells.Item(122,10).Value = 0
$ws.Cells.Item(122,11).Value = 9114.500100000001
$ws.Cells.Item(122,12).Value = 0
$ws.Cells.Item(122,13).Value = -6664.500100000001
$ws.Cells.Item(122,14).ClearContents()
# Row 127
$ws.Cells.Item(127,8).Value = 109897.164
$ws.Cells.Item(127,10).Value = 111876.4
$ws.Cells.Item(127,12).Value = 111876.4
$ws.Cells.Item(127,14).Value = -121796.4
# Row 132
$ws.Cells.Item(132,8).Value = 2344.2942
$ws.Cells.Item(132,9).Value = 2346.6428
$ws.Cells.Item(132,11).Value = 7039.928400000001
$ws.Cells.Item(132,13).Value = -4509.928400000001
# Row 141
$ws.Cells.Item(141,8).Value = 561723.4
$ws.Cells.Item(141,10).Value = 561723.4
$ws.Cells.Item(141,12).Value = 561723.4
$ws.Cells.Item(141,14).Value = -572083.4

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Cells.Item(2,8).Value = 60
$ws.Cells.Item(2,9).Value = 80
$ws.Cells.Item(2,10).Value = 20
$ws.Cells.Item(2,11).Value = 480
$ws.Cells.Item(2,12).Value = 120
$ws.Cells.Item(2,13).Value = -367
$ws.Cells.Item(2,14).Value = -346
# Row 3
$ws.Cells.Item(3,8).Value = 7623.25
$ws.Cells.Item(3,9).Value = 7279
$ws.Cells.Item(3,11).Value = 21837
$ws.Cells.Item(3,13).Value = -21725
# Row 34
$ws.Cells.Item(34,8).Value = 881.36365
$ws.Cells.Item(34,10).Value = 1800
$ws.Cells.Item(34,12).Value = 5400
$ws.Cells.Item(34,14).Value = -5568
# Row 39
$ws.Cells.Item(39,8).Value = 4333.3335
$ws.Cells.Item(39,10).Value = 5000
$ws.Cells.Item(39,12).Value = 15000
$ws.Cells.Item(39,14).Value = -15588
# Row 139
$ws.Cells.Item(139,8).Value = 4244
$ws.Cells.Item(139,9).Value = 3956.5715
$ws.Cells.Item(139,11).Value = 11869.7145
$ws.Cells.Item(139,13).Value = -6729.7145

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 32
$ws.Cells.Item(32,8).Value = 59949.25
$ws.Cells.Item(32,10).Value = 59949.25
$ws.Cells.Item(32,12).Value = 59949.25
$ws.Cells.Item(32,14).Value = -60541.25
# Row 45
$ws.Cells.Item(45,8).Value = 103326
$ws.Cells.Item(45,10).Value = 103326
$ws.Cells.Item(45,12).Value = 103326
$ws.Cells.Item(45,14).Value = -104444
# Row 64
$ws.Cells.Item(64,8).Value = 0
$ws.Cells.Item(64,10).Value = 0
$ws.Cells.Item(64,12).Value = 0
$ws.Cells.Item(64,14).ClearContents()
# Row 67
$ws.Cells.Item(67,8).Value = 0
$ws.Cells.Item(67,10).Value = 0
$ws.Cells.Item(67,12).Value = 0
$ws.Cells.Item(67,14).ClearContents()
# Row 97
$ws.Cells.Item(97,8).Value = 1246.5
$ws.Cells.Item(97,9).Value = 895
$ws.Cells.Item(97,10).Value = 1949.5
$ws.Cells.Item(97,11).Value = 895
$ws.Cells.Item(97,12).Value = 1949.5
$ws.Cells.Item(97,13).Value = -399
$ws.Cells.Item(97,14).Value = -2941.5
# Row 122
$ws.Cells.Item(122,8).Value = 2490.95
$ws.Cells.Item(122,9).Value = 2771.1428
$ws.Cells.Item(122,10).Value = 1837.1666
$ws.Cells.Item(122,11).Value = 8313.428400000001
$ws.Cells.Item(122,12).Value = 5511.4998
$ws.Cells.Item(122,13).Value = -5863.428400000001
$ws.Cells.Item(122,14).Value = -10411.4998
# Row 131
$ws.Cells.Item(131,8).Value = 44081.5
$ws.Cells.Item(131,10).Value = 44081.5
$ws.Cells.Item(131,12).Value = 44081.5
$ws.Cells.Item(131,14).Value = -54161.5
# Row 132
$ws.Cells.Item(132,8).Value = 4182.8667
$ws.Cells.Item(132,9).Value = 4195.9287
$ws.Cells.Item(132,10).Value = 4000
$ws.Cells.Item(132,11).Value = 12587.7861
$ws.Cells.Item(132,12).Value = 12000
$ws.Cells.Item(132,13).Value = -10057.7861
$ws.Cells.Item(132,14).Value = -17060
# Row 136
$ws.Cells.Item(136,8).Value = 28716.424
$ws.Cells.Item(136,10).Value = 28716.424
$ws.Cells.Item(136,12).Value = 86149.272
$ws.Cells.Item(136,14).Value = -91249.272

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40,8).Value = 7460.25
$ws.Cells.Item(40,9).Value = 7837.9375
$ws.Cells.Item(40,11).Value = 7837.9375
$ws.Cells.Item(40,13).Value = -7701.9375
# Row 88
$ws.Cells.Item(88,8).Value = 52421.855
$ws.Cells.Item(88,9).Value = 53658.832
$ws.Cells.Item(88,11).Value = 53658.832
$ws.Cells.Item(88,13).Value = -53230.832
# Row 91
$ws.Cells.Item(91,8).Value = 52421.855
$ws.Cells.Item(91,9).Value = 53658.832
$ws.Cells.Item(91,11).Value = 53658.832
$ws.Cells.Item(91,13).Value = -52176.832
# Row 122
$ws.Cells.Item(122,8).Value = 3414.923
$ws.Cells.Item(122,9).Value = 3429.5
$ws.Cells.Item(122,10).Value = 3391.6
$ws.Cells.Item(122,11).Value = 10288.5
$ws.Cells.Item(122,12).Value = 10174.8
$ws.Cells.Item(122,13).Value = -7838.5
$ws.Cells.Item(122,14).Value = -15074.8
# Row 123
$ws.Cells.Item(123,8).Value = 0
$ws.Cells.Item(123,10).Value = 0
$ws.Cells.Item(123,12).Value = 0
$ws.Cells.Item(123,14).ClearContents()
# Row 136
$ws.Cells.Item(136,8).Value = 3004.675
$ws.Cells.Item(136,9).Value = 2225.6785
$ws.Cells.Item(136,11).Value = 6677.0355
$ws.Cells.Item(136,13).Value = -4127.0355

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Cells.Item(81,8).Value = 2972.1538
$ws.Cells.Item(81,10).Value = 3966.1667
$ws.Cells.Item(81,12).Value = 7932.3334
$ws.Cells.Item(81,14).Value = -10054.3334
# Row 84
$ws.Cells.Item(84,8).Value = 2972.1538
$ws.Cells.Item(84,10).Value = 3966.1667
$ws.Cells.Item(84,12).Value = 39661.667
$ws.Cells.Item(84,14).Value = -50269.667
# Row 100
$ws.Cells.Item(100,8).Value = 5377.8
$ws.Cells.Item(100,9).Value = 5629.6665
$ws.Cells.Item(100,10).Value = 5000
$ws.Cells.Item(100,11).Value = 11259.333
$ws.Cells.Item(100,12).Value = 10000
$ws.Cells.Item(100,13).Value = -10718.333
$ws.Cells.Item(100,14).Value = -11082
# Row 113
$ws.Cells.Item(113,8).Value = 430.6154
$ws.Cells.Item(113,9).Value = 530.6667
$ws.Cells.Item(113,10).Value = 344.85715
$ws.Cells.Item(113,11).Value = 1592.0001
$ws.Cells.Item(113,12).Value = 1034.57145
$ws.Cells.Item(113,13).Value = 577.9999
$ws.Cells.Item(113,14).Value = -5374.571449999999
# Row 132
$ws.Cells.Item(132,8).Value = 1526.0286
$ws.Cells.Item(132,9).Value = 1356.4375
$ws.Cells.Item(132,10).Value = 3335
$ws.Cells.Item(132,11).Value = 4069.3125
$ws.Cells.Item(132,12).Value = 10005
$ws.Cells.Item(132,13).Value = -1539.3125
$ws.Cells.Item(132,14).Value = -15065
# Row 136
$ws.Cells.Item(136,8).Value = 22743.1
$ws.Cells.Item(136,9).Value = 1315.2693
$ws.Cells.Item(136,10).Value = 45956.582
$ws.Cells.Item(136,11).Value = 3945.8079
$ws.Cells.Item(136,12).Value = 137869.746
$ws.Cells.Item(136,13).Value = -1395.8079
$ws.Cells.Item(136,14).Value = -142969.746
